$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 113
$ws1.Range("F4").Value  = 1570
$ws1.Range("F5").Value  = 264
$ws1.Range("F6").Value  = 59
$ws1.Range("F7").Value  = 1323
$ws1.Range("F8").Value  = 10187
$ws1.Range("F10").Value = 136
$ws1.Range("F13").Value = 390
$ws1.Range("F14").Value = 7058
$ws1.Range("F15").Value = 1099
$ws1.Range("F16").Value = 659
$ws1.Range("F17").Value = 28
$ws1.Range("F19").Value = 230

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 8
$ws2.Range("F3").Value = 555

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 113
$ws4.Range("F4").Value  = 1570
$ws4.Range("F5").Value  = 264
$ws4.Range("F6").Value  = 8
$ws4.Range("F7").Value  = 59
$ws4.Range("F8").Value  = 1323
$ws4.Range("F9").Value  = 555
$ws4.Range("F11").Value = 10187
$ws4.Range("F13").Value = 136
$ws4.Range("F16").Value = 390
$ws4.Range("F17").Value = 7058
$ws4.Range("F18").Value = 1099
$ws4.Range("F19").Value = 659
$ws4.Range("F20").Value = 28
$ws4.Range("F22").Value = 230
